$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new survey response was recorded (row 13). Clone the formatting of the
# last existing data row (12) down onto the new row, then fill in its values.
$ws.Range("A12:R12").Copy()
$ws.Range("A13:R13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(13).RowHeight = 39.95

$ws.Range("A13").Value = "Business"
$ws.Range("B13").Value = "Female"
$ws.Range("C13").Value = 89
$ws.Range("D13").Value = "Region 9"
$ws.Range("E13").Value = "ASDS"
$ws.Range("F13").Value = "Other request or inqueries"
$ws.Range("G13").Value = "I do not know what a CC is and I did not see on in this office."
$ws.Range("H13").Value = "N/A"
$ws.Range("I13").Value = "N/A"
$ws.Range("J13:R13").Value = "strongly-agree"
